$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 202.7037
$ws.Range("I53").Value = 206.14285
$ws.Range("J53").Value = 199
$ws.Range("K53").Value = 206.14285
$ws.Range("L53").Value = 199
$ws.Range("M53").Value = 430.85715
$ws.Range("N53").Value = -1473
$ws.Range("H80").Value = 1251.7778
$ws.Range("I80").Value = 1750
$ws.Range("J80").Value = 1109.4286
$ws.Range("K80").Value = 5250
$ws.Range("L80").Value = 3328.2858
$ws.Range("M80").Value = -4252
$ws.Range("N80").Value = -5324.2858
$ws.Range("H83").Value = 1251.7778
$ws.Range("I83").Value = 1750
$ws.Range("J83").Value = 1109.4286
$ws.Range("K83").Value = 15750
$ws.Range("L83").Value = 9984.857399999999
$ws.Range("M83").Value = -10758
$ws.Range("N83").Value = -19968.8574
$ws.Range("H86").Value = 2067.3635
$ws.Range("I86").Value = 2155.111
$ws.Range("J86").Value = 1672.5
$ws.Range("K86").Value = 2155.111
$ws.Range("L86").Value = 1672.5
$ws.Range("M86").Value = -1032.111
$ws.Range("N86").Value = -3918.5
$ws.Range("H89").Value = 2067.3635
$ws.Range("I89").Value = 2155.111
$ws.Range("J89").Value = 1672.5
$ws.Range("K89").Value = 10775.555
$ws.Range("L89").Value = 8362.5
$ws.Range("M89").Value = -5159.555
$ws.Range("N89").Value = -19594.5
$ws.Range("H113").Value = 4492.231
$ws.Range("I113").Value = 3849.8333
$ws.Range("J113").Value = 5042.857
$ws.Range("K113").Value = 3849.8333
$ws.Range("L113").Value = 5042.857
$ws.Range("M113").Value = -595.8332999999998
$ws.Range("N113").Value = -11550.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 725630.1
$ws.Range("I45").Value = 2527129.8
$ws.Range("J45").Value = 5030.3
$ws.Range("K45").Value = 2527129.8
$ws.Range("L45").Value = 5030.3
$ws.Range("M45").Value = -2526752.8
$ws.Range("N45").Value = -5784.3
$ws.Range("H132").Value = 5794.408
$ws.Range("I132").Value = 6420.641
$ws.Range("J132").Value = 3352.1
$ws.Range("K132").Value = 19261.923
$ws.Range("L132").Value = 10056.3
$ws.Range("M132").Value = -16731.923
$ws.Range("N132").Value = -15116.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2020.5667
$ws.Range("I86").Value = 1866.5264
$ws.Range("J86").Value = 2286.6365
$ws.Range("K86").Value = 1866.5264
$ws.Range("L86").Value = 2286.6365
$ws.Range("M86").Value = -743.5264
$ws.Range("N86").Value = -4532.636500000001
$ws.Range("H89").Value = 2020.5667
$ws.Range("I89").Value = 1866.5264
$ws.Range("J89").Value = 2286.6365
$ws.Range("K89").Value = 9332.632
$ws.Range("L89").Value = 11433.1825
$ws.Range("M89").Value = -3716.632
$ws.Range("N89").Value = -22665.1825
$ws.Range("H107").Value = 13228.467
$ws.Range("I107").Value = 3686.6924
$ws.Range("J107").Value = 75250
$ws.Range("K107").Value = 3686.6924
$ws.Range("L107").Value = 75250
$ws.Range("M107").Value = -1766.6924
$ws.Range("N107").Value = -79090
$ws.Range("H134").Value = 2897.4
$ws.Range("I134").Value = 2562.6667
$ws.Range("J134").Value = 3399.5
$ws.Range("K134").Value = 7688.000100000001
$ws.Range("L134").Value = 10198.5
$ws.Range("M134").Value = -5153.000100000001
$ws.Range("N134").Value = -15268.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1837.75
$ws.Range("I16").Value = 1784.4
$ws.Range("J16").Value = 1926.6666
$ws.Range("K16").Value = 1784.4
$ws.Range("L16").Value = 1926.6666
$ws.Range("M16").Value = -1497.4
$ws.Range("N16").Value = -2500.6666
$ws.Range("H31").Value = 2418.658
$ws.Range("I31").Value = 1848.1052
$ws.Range("J31").Value = 2989.2104
$ws.Range("K31").Value = 1848.1052
$ws.Range("L31").Value = 2989.2104
$ws.Range("M31").Value = -1553.1052
$ws.Range("N31").Value = -3579.2104
$ws.Range("H34").Value = 2418.658
$ws.Range("I34").Value = 1848.1052
$ws.Range("J34").Value = 2989.2104
$ws.Range("K34").Value = 1848.1052
$ws.Range("L34").Value = 2989.2104
$ws.Range("M34").Value = -1646.1052
$ws.Range("N34").Value = -3393.2104
$ws.Range("H87").Value = 25330
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 25330
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 25330
$ws.Range("N87").Value = -27702
$ws.Range("H90").Value = 25330
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 25330
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 75990
$ws.Range("N90").Value = -87846
$ws.Range("H99").Value = 7721.875
$ws.Range("I99").Value = 1724.3334
$ws.Range("J99").Value = 25714.5
$ws.Range("K99").Value = 1724.3334
$ws.Range("L99").Value = 25714.5
$ws.Range("M99").Value = -226.3334
$ws.Range("N99").Value = -28710.5
$ws.Range("H107").Value = 1378.3
$ws.Range("I107").Value = 910.1111
$ws.Range("J107").Value = 1761.3636
$ws.Range("K107").Value = 910.1111
$ws.Range("L107").Value = 1761.3636
$ws.Range("M107").Value = 1009.8889
$ws.Range("N107").Value = -5601.3636
$ws.Range("H113").Value = 1837.75
$ws.Range("I113").Value = 1784.4
$ws.Range("J113").Value = 1926.6666
$ws.Range("K113").Value = 1784.4
$ws.Range("L113").Value = 1926.6666
$ws.Range("M113").Value = 385.5999999999999
$ws.Range("N113").Value = -6266.6666
$ws.Range("H126").Value = 7721.875
$ws.Range("I126").Value = 1724.3334
$ws.Range("J126").Value = 25714.5
$ws.Range("K126").Value = 5173.0002
$ws.Range("L126").Value = 77143.5
$ws.Range("M126").Value = -2703.0002
$ws.Range("N126").Value = -82083.5
$ws.Range("H134").Value = 2395.838
$ws.Range("I134").Value = 2120.2122
$ws.Range("J134").Value = 4669.75
$ws.Range("K134").Value = 6360.6366
$ws.Range("L134").Value = 14009.25
$ws.Range("M134").Value = -3825.6366
$ws.Range("N134").Value = -19079.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1522.4445
$ws.Range("I107").Value = 606.7143
$ws.Range("J107").Value = 1842.95
$ws.Range("K107").Value = 1820.1429
$ws.Range("L107").Value = 5528.85
$ws.Range("M107").Value = 99.85710000000017
$ws.Range("N107").Value = -9368.85
$ws.Range("H136").Value = 4411.4287
$ws.Range("I136").Value = 1612.3334
$ws.Range("J136").Value = 9449.799999999999
$ws.Range("K136").Value = 4837.0002
$ws.Range("L136").Value = 28349.4
$ws.Range("M136").Value = 262.9997999999996
$ws.Range("N136").Value = -38549.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2749.5
$ws.Range("J22").Value = 2333.6667
$ws.Range("K22").Value = 2749.5
$ws.Range("L22").Value = 2333.6667
$ws.Range("M22").Value = -2454.5
$ws.Range("N22").Value = -2923.6667
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2749.5
$ws.Range("J27").Value = 2333.6667
$ws.Range("K27").Value = 2749.5
$ws.Range("L27").Value = 2333.6667
$ws.Range("M27").Value = -2642.5
$ws.Range("N27").Value = -2547.6667
$ws.Range("H61").Value = 207461.72
$ws.Range("I61").Value = 277749.75
$ws.Range("J61").Value = 113744.336
$ws.Range("K61").Value = 277749.75
$ws.Range("L61").Value = 113744.336
$ws.Range("M61").Value = -277547.75
$ws.Range("N61").Value = -114148.336
$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 3500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3500
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2751
$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 3500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 17500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -13756
$ws.Range("H93").Value = 25681.785
$ws.Range("I93").Value = 2095.4285
$ws.Range("J93").Value = 49268.145
$ws.Range("K93").Value = 2095.4285
$ws.Range("L93").Value = 49268.145
$ws.Range("M93").Value = -847.4285
$ws.Range("N93").Value = -51764.145
$ws.Range("H107").Value = 2061.5
$ws.Range("I107").Value = 2061.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2061.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -141.5
$ws.Range("H113").Value = 207461.72
$ws.Range("I113").Value = 277749.75
$ws.Range("J113").Value = 113744.336
$ws.Range("K113").Value = 277749.75
$ws.Range("L113").Value = 113744.336
$ws.Range("M113").Value = -275579.75
$ws.Range("N113").Value = -118084.336
$ws.Range("H132").Value = 9404.134
$ws.Range("I132").Value = 14796.5
$ws.Range("J132").Value = 3241.4285
$ws.Range("K132").Value = 44389.5
$ws.Range("L132").Value = 9724.2855
$ws.Range("M132").Value = -41859.5
$ws.Range("N132").Value = -14784.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1909.8235
$ws.Range("I81").Value = 2039
$ws.Range("J81").Value = 941
$ws.Range("K81").Value = 4078
$ws.Range("L81").Value = 1882
$ws.Range("M81").Value = -3017
$ws.Range("N81").Value = -4004
$ws.Range("H84").Value = 1909.8235
$ws.Range("I84").Value = 2039
$ws.Range("J84").Value = 941
$ws.Range("K84").Value = 20390
$ws.Range("L84").Value = 9410
$ws.Range("M84").Value = -15086
$ws.Range("N84").Value = -20018
$ws.Range("H107").Value = 35715916
$ws.Range("I107").Value = 1795.1
$ws.Range("J107").Value = 125001220
$ws.Range("K107").Value = 5385.299999999999
$ws.Range("L107").Value = 375003660
$ws.Range("M107").Value = -3465.299999999999
$ws.Range("N107").Value = -375007500
$ws.Range("H132").Value = 8430.686
$ws.Range("I132").Value = 8511.333000000001
$ws.Range("J132").Value = 7100
$ws.Range("K132").Value = 25533.999
$ws.Range("L132").Value = 21300
$ws.Range("M132").Value = -23003.999
$ws.Range("N132").Value = -26360
